$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 63-65: results for these three fixtures were re-matched to
#    the correct games. Columns F:V (everything except the shared
#    Indice/pais/torneio/temporada/data_partida columns A:E) rotate:
#       new row63 <- old row64
#      new row64 <- old row65
#      new row65 <- old row63
# ------------------------------------------------------------------

$oldRow63 = @{
    F = $ws.Range("F63").Value()
    G = $ws.Range("G63").Value()
    H = $ws.Range("H63").Value()
    I = $ws.Range("I63").Value()
    J = $ws.Range("J63").Value()
    K = $ws.Range("K63").Value()
    L = $ws.Range("L63").Value()
    M = $ws.Range("M63").Value()
    N = $ws.Range("N63").Value()
    O = $ws.Range("O63").Value()
    P = $ws.Range("P63").Value()
    Q = $ws.Range("Q63").Value()
    R = $ws.Range("R63").Value()
    S = $ws.Range("S63").Value()
    T = $ws.Range("T63").Value()
    U = $ws.Range("U63").Value()
    V = $ws.Range("V63").Value()
}

$oldRow64 = @{
    F = $ws.Range("F64").Value()
    G = $ws.Range("G64").Value()
    H = $ws.Range("H64").Value()
    I = $ws.Range("I64").Value()
    J = $ws.Range("J64").Value()
    K = $ws.Range("K64").Value()
    L = $ws.Range("L64").Value()
    M = $ws.Range("M64").Value()
    N = $ws.Range("N64").Value()
    O = $ws.Range("O64").Value()
    P = $ws.Range("P64").Value()
    Q = $ws.Range("Q64").Value()
    R = $ws.Range("R64").Value()
    S = $ws.Range("S64").Value()
    T = $ws.Range("T64").Value()
    U = $ws.Range("U64").Value()
    V = $ws.Range("V64").Value()
}

$oldRow65 = @{
    F = $ws.Range("F65").Value()
    G = $ws.Range("G65").Value()
    H = $ws.Range("H65").Value()
    I = $ws.Range("I65").Value()
    J = $ws.Range("J65").Value()
    K = $ws.Range("K65").Value()
    L = $ws.Range("L65").Value()
    M = $ws.Range("M65").Value()
    N = $ws.Range("N65").Value()
    O = $ws.Range("O65").Value()
    P = $ws.Range("P65").Value()
    Q = $ws.Range("Q65").Value()
    R = $ws.Range("R65").Value()
    S = $ws.Range("S65").Value()
    T = $ws.Range("T65").Value()
    U = $ws.Range("U65").Value()
    V = $ws.Range("V65").Value()
}

foreach ($col in @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")) {
    $ws.Range("$col" + "63").Value = $oldRow64[$col]
    $ws.Range("$col" + "64").Value = $oldRow65[$col]
    $ws.Range("$col" + "65").Value = $oldRow63[$col]
}

# ------------------------------------------------------------------
# 2) Append 4 new match rows (112-115), continuing the Indice sequence
#    and formatted the same way as the existing data rows.
# ------------------------------------------------------------------

$ws.Range("A111:V111").Copy()
$ws.Range("A112:V115").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 112; A = 111; B = "spain"; C = "laliga"; D = "2023-2024"; E = 45234.58333333334;
       F = "Osasuna"; G = 2; H = "Girona"; I = 4;
       J = 2.39; K = "22/10/2023 12:02"; L = 2.99; M = "04/11/2023 13:59";
       N = 3.31; O = "22/10/2023 12:02"; P = 3.31; Q = "04/11/2023 13:58";
       R = 3.19; S = "22/10/2023 12:02"; T = 2.56; U = "04/11/2023 13:59";
       V = "https://www.betexplorer.com/football/spain/laliga/osasuna-girona/GncXCyE5/" },
    @{ Row = 113; A = 112; B = "spain"; C = "laliga"; D = "2023-2024"; E = 45234.67708333334;
       F = "Betis"; G = 2; H = "Mallorca"; I = 0;
       J = 1.97; K = "22/10/2023 12:02"; L = 1.9; M = "04/11/2023 16:13";
       N = 3.28; O = "22/10/2023 12:02"; P = 3.45; Q = "04/11/2023 16:13";
       R = 4.51; S = "22/10/2023 12:02"; T = 4.68; U = "04/11/2023 16:13";
       V = "https://www.betexplorer.com/football/spain/laliga/betis-mallorca/IiR2cdEI/" },
    @{ Row = 114; A = 113; B = "spain"; C = "laliga"; D = "2023-2024"; E = 45234.77083333334;
       F = "Celta Vigo"; G = 1; H = "Sevilla"; I = 1;
       J = 2.13; K = "22/10/2023 12:02"; L = 2.41; M = "04/11/2023 18:24";
       N = 3.4; O = "22/10/2023 12:02"; P = 3.41; Q = "04/11/2023 18:29";
       R = 3.43; S = "22/10/2023 12:02"; T = 3.13; U = "04/11/2023 18:22";
       V = "https://www.betexplorer.com/football/spain/laliga/celta-vigo-sevilla/drQ6dGTO/" },
    @{ Row = 115; A = 114; B = "spain"; C = "laliga"; D = "2023-2024"; E = 45234.875;
       F = "Real Sociedad"; G = 0; H = "Barcelona"; I = 1;
       J = 3.21; K = "22/10/2023 12:02"; L = 3.38; M = "04/11/2023 20:57";
       N = 3.4; O = "22/10/2023 12:02"; P = 3.29; Q = "04/11/2023 20:57";
       R = 2.22; S = "22/10/2023 12:02"; T = 2.33; U = "04/11/2023 20:59";
       V = "https://www.betexplorer.com/football/spain/laliga/real-sociedad-barcelona/rZbPEcqg/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
}
